# Update countries & provincias Spain
# - Arabia Saudita overtakes Belgica (rows 18/19 swap country + stats)
# - Republica de Macedonia overtakes Islandia (rows 86/87 swap country + stats)
# - Tayikistan overtakes Lituania (rows 89/90 swap country + stats)
# - Refresh case counts for Alemania (11), India (14), Senegal (79)
# - Bump the "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18/19: Belgica / Arabia Saudita swap rank, with refreshed stats ---
$ws.Range("A18").Value2 = "Arabia Saudita"
$ws.Range("B18").Value2 = 57345
$ws.Range("C18").Value2 = 2593
$ws.Range("D18").Value2 = 28748
$ws.Range("E18").Value2 = 28277
$ws.Range("F18").Value2 = 0
$ws.Range("G18").Value2 = 8
$ws.Range("H18").Value2 = 320

$ws.Range("A19").Value2 = "Belgica"
$ws.Range("B19").Value2 = 55559
$ws.Range("C19").Value2 = 279
$ws.Range("D19").Value2 = 14657
$ws.Range("E19").Value2 = 31822
$ws.Range("F19").Value2 = 0
$ws.Range("G19").Value2 = 28
$ws.Range("H19").Value2 = 9080

# --- Row 11: Alemania stats refresh ---
$ws.Range("B11").Value2 = 176807
$ws.Range("C11").Value2 = 156
$ws.Range("D11").Value2 = 154600
$ws.Range("E11").Value2 = 14153
$ws.Range("F11").Value2 = 0
$ws.Range("G11").Value2 = 5
$ws.Range("H11").Value2 = 8054

# --- Row 14: India stats refresh ---
$ws.Range("B14").Value2 = 96492
$ws.Range("C14").Value2 = 794
$ws.Range("D14").Value2 = 36824
$ws.Range("E14").Value2 = 56627
$ws.Range("F14").Value2 = 0
$ws.Range("G14").Value2 = 16
$ws.Range("H14").Value2 = 3041

# --- Row 79: Senegal stats refresh ---
$ws.Range("B79").Value2 = 2544
$ws.Range("C79").Value2 = 64
$ws.Range("D79").Value2 = 1076
$ws.Range("E79").Value2 = 1442
$ws.Range("F79").Value2 = 0
$ws.Range("G79").Value2 = 1
$ws.Range("H79").Value2 = 26

# --- Row 86/87: Islandia / Republica de Macedonia swap rank, with refreshed stats ---
$ws.Range("A86").Value2 = "Republica de Macedonia"
$ws.Range("B86").Value2 = 1817
$ws.Range("C86").Value2 = 25
$ws.Range("D86").Value2 = 1301
$ws.Range("E86").Value2 = 412
$ws.Range("F86").Value2 = 0
$ws.Range("G86").Value2 = 3
$ws.Range("H86").Value2 = 104

$ws.Range("A87").Value2 = "Islandia"
$ws.Range("B87").Value2 = 1802
$ws.Range("C87").Value2 = 0
$ws.Range("D87").Value2 = 1786
$ws.Range("E87").Value2 = 6
$ws.Range("F87").Value2 = 0
$ws.Range("G87").Value2 = 0
$ws.Range("H87").Value2 = 10

# --- Row 89/90: Lituania / Tayikistan swap rank, with refreshed stats ---
$ws.Range("A89").Value2 = "Tayikistan"
$ws.Range("B89").Value2 = 1729
$ws.Range("C89").Value2 = 205
$ws.Range("D89").Value2 = 0
$ws.Range("E89").Value2 = 1688
$ws.Range("F89").Value2 = 0
$ws.Range("G89").Value2 = 2
$ws.Range("H89").Value2 = 41

$ws.Range("A90").Value2 = "Lituania"
$ws.Range("B90").Value2 = 1547
$ws.Range("C90").Value2 = 6
$ws.Range("D90").Value2 = 997
$ws.Range("E90").Value2 = 491
$ws.Range("F90").Value2 = 0
$ws.Range("G90").Value2 = 3
$ws.Range("H90").Value2 = 59

# --- Timestamp refresh ---
$ws.Range("A1").Value2 = "Datos actualizados a 18 de Mayo de 2020 a las 15:05"
